$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" Date placeholder text from
#    5/2/21 -> 5/4/21 everywhere it appears: the slide master and every
#    slide layout (PlaceholderFormat.Type 16 = ppPlaceholderDate).
# ---------------------------------------------------------------------------
$newDate = "5/4/21"

function Update-DatePlaceholder($shapes) {
  for ($i = 1; $i -le $shapes.Count; $i++) {
    $sh = $shapes.Item($i)
    if ($sh.Type -eq 14 -and $sh.HasTextFrame) {
      $ph = $sh.PlaceholderFormat
      if ($ph.Type -eq 16) {
        $sh.TextFrame.TextRange.Text = $newDate
      }
    }
  }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout hanging off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
  $layout = $layouts.Item($li)
  Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Rename the "CFE Features" label to "NL Features" and shrink/reposition
#    its textbox to match the new (shorter) auto-fit text box size.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
  $sh = $shapes.Item($i)
  if ($sh.HasTextFrame) {
    if ($sh.TextFrame.TextRange.Text -eq "CFE Features") {
      $sh.Left = 2058333 / 12700
      $sh.Width = 678391 / 12700
      $sh.TextFrame.TextRange.Text = "NL Features"
    }
  }
}
